# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# Before:  2021-Q2 | 2021-Q3 | 2021-Q4 | 总计
# After :  2021-Q2 | 2021-Q3 | 2021-Q4 | 2022-Q1 | 总计
#
# The old "总计" worksheet (sheetId 4) is repurposed in-place to become the
# new "2022-Q1" fund-holdings sheet (keeps sheetId/rId 4), and a brand new
# "总计" worksheet (sheetId/rId 5) is appended at the end, carrying the
# former totals table plus a freshly prepended 2022-Q1 summary row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$fundTemplate = $wb.Worksheets.Item(3)     # "2021-Q4" - 8-column fund holdings layout
$oldTotal     = $wb.Worksheets.Item($wb.Worksheets.Count)   # "总计" (sheetId 4)

# 1) Duplicate the existing "总计" sheet; the duplicate lands right after it
#    and inherits the next free sheetId (5) together with a fresh rId.
$oldTotal.Copy($null, $oldTotal)

$newQ1    = $wb.Worksheets.Item(4)   # formerly "总计" (sheetId 4) -> becomes "2022-Q1"
$newTotal = $wb.Worksheets.Item(5)   # the duplicate (sheetId 5)   -> becomes "总计"

$newQ1.Name    = "2022-Q1"
$newTotal.Name = "总计"

# ---------------------------------------------------------------------------
# 2) Rebuild "2022-Q1" as an 8-column fund-holdings sheet (A1:H4), matching
#    the layout/style used by the quarterly sheets.
# ---------------------------------------------------------------------------

# Pull header formatting (incl. the bold/bordered style) from the template.
$fundTemplate.Range("B1:H1").Copy($newQ1.Range("B1:H1"))
$fundTemplate.Range("A2:A4").Copy($newQ1.Range("A2:A4"))

$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

$newQ1.Range("A2").Value = 0
$newQ1.Range("A3").Value = 1
$newQ1.Range("A4").Value = 2

# Columns B:G hold text values (fund codes/names/figures stored as strings in
# the source data) - force text storage so numeric-looking strings such as
# "23.39" are not silently coerced into numbers.
$newQ1.Range("B2:G4").NumberFormat = "@"

$newQ1.Range("B2").Value = "519692"
$newQ1.Range("C2").Value = "交银成长混合A"
$newQ1.Range("D2").Value = "23.39"
$newQ1.Range("E2").Value = "82.72"
$newQ1.Range("F2").Value = "2.65"
$newQ1.Range("G2").Value = "0.6198"
$newQ1.Range("H2").Value = 10

$newQ1.Range("B3").Value = "960016"
$newQ1.Range("C3").Value = "交银成长混合H"
$newQ1.Range("D3").Value = "23.39"
$newQ1.Range("E3").Value = "82.72"
$newQ1.Range("F3").Value = "2.65"
$newQ1.Range("G3").Value = "0.6198"
$newQ1.Range("H3").Value = 10

$newQ1.Range("B4").Value = "002567"
$newQ1.Range("C4").Value = "大成国家安全主题灵活配置混合"
$newQ1.Range("D4").Value = "0.34"
$newQ1.Range("E4").Value = "52.90"
$newQ1.Range("F4").Value = "3.21"
$newQ1.Range("G4").Value = "0.0109"
$newQ1.Range("H4").Value = 10

# Restore the default cell style on the data cells (only NumberFormat should
# stick around as the distinguishing trait, matching the original workbook's
# styling where data rows carry no explicit style index).
$newQ1.Range("B2:G4").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Rebuild "总计" with the 2022-Q1 row prepended to the historical table.
# ---------------------------------------------------------------------------

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

# The duplicated sheet only had 4 rows (header + 3 data); stamp the same
# index-column style onto the newly-needed 5th row before filling it in.
$newTotal.Range("A4").Copy($newTotal.Range("A5"))

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 3
$newTotal.Range("D2").Value = 1.25

$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 4
$newTotal.Range("D3").Value = 0.26

$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 3
$newTotal.Range("D4").Value = 0.45

$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 2
$newTotal.Range("D5").Value = 0.05
